# FA120_TestData_ManuallyManageFinancialsTransactions_21C.xlsx
# "Add files via upload" - scrub the hard-coded Oracle Cloud URL / username /
# password that used to live in X2:Z2 (and the live hyperlink that pointed
# at the URL cell) on the Input_Value sheet, and update the view/selection
# state left behind by the author's last interactive session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")
$ws.Activate()

# Remove the hyperlink that lived on X2 (https://edrx.fa.us2.oraclecloud.com/)
# and blank out the credential cells it (and its neighbours) held.
$ws.Range("X2").Hyperlinks.Delete()
$ws.Range("X2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()

# Leave the sheet scrolled/selected over the now-empty X2:Z2 block, matching
# the view state captured the last time the workbook was saved.
$ws.Range("X2:Z2").Select()
